$wb = $excel.ActiveWorkbook

# Add the new "News" worksheet after the last existing sheet (LoginPage)
$lastSheet = $wb.Worksheets.Item($wb.Worksheets.Count)
$news = $wb.Worksheets.Add($null, $lastSheet)
$news.Name = "News"

# Populate cell A1 with the news text
$news.Range("A1").Value = "This is a news written in excel sheet"

# Make the News sheet the active sheet/tab
$news.Activate()
